# Cập nhật tên bài tập
# 1) Slide 1 (title slide) subtitle: merge "Bài " + "23. " + "Thu " runs
#    into a single run "Bài 23. Thu " (same text, now one run).
$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$subtitleRange = $subtitle.TextFrame.TextRange
$subtitleChars = $subtitleRange.Characters(1, 12)
$subtitleChars.Text = "Bài 23. Thu "

# 2) Slide 28 title "Bài tập" -> "Bài " + "tập 23.1"
$s28 = $p.Slides.Item(28)
$title28 = $s28.Shapes.Item(2)
$title28Range = $title28.TextFrame.TextRange
$title28Chars = $title28Range.Characters(5, 3)
$title28Chars.Text = "tập 23.1"
